$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "71.540.00"
$ws.Range("E2").Value = "  -1.74%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.878.05"
$ws.Range("E3").Value = "  -2.64%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.06%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "603.97"
$ws.Range("E5").Value = "  -1.78%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "173.11"
$ws.Range("E6").Value = "  +5.51%  "

$ws.Range("E7").Value = "  -2.15%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.999"
$ws.Range("E8").Value = "  +0.07%  "

$ws.Range("E9").Value = "  -0.73%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.178"
$ws.Range("E10").Value = "  +5.61%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "53.97"
$ws.Range("E11").Value = "  -0.53%  "

$ws.Range("E12").Value = "  +0.88%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "11.55"
$ws.Range("E13").Value = "  +4.86%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.493.60"
$ws.Range("E14").Value = "  -2.68%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "21.17"
$ws.Range("E15").Value = "  +3.05%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.874.52"
$ws.Range("E16").Value = "  -2.89%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "13.99"
$ws.Range("E17").Value = "  -1.15%  "

$ws.Range("E18").Value = "  -3.64%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "71.257.21"
$ws.Range("E20").Value = "  -1.73%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "441.13"
$ws.Range("E21").Value = "  +0.27%  "

$ws.Range("E22").Value = "  -1.72%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "94.64"
$ws.Range("E23").Value = "  -2.16%  "

$ws.Range("E24").Value = "  -4.57%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "13.89"
$ws.Range("E25").Value = "  -3.64%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "11.71"
$ws.Range("E26").Value = "  +1.92%  "

$ws.Range("E27").Value = "  -6.03%  "

$ws.Range("E28").Value = "  +0.37%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "10.47"
$ws.Range("E29").Value = "  -0.97%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.58"
$ws.Range("E30").Value = "  +10.21%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "35.28"
$ws.Range("E31").Value = "  -3.33%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "13.58"
$ws.Range("E32").Value = "  -2.89%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "47.94"
$ws.Range("E33").Value = "  -1.78%  "

$ws.Range("E34").Value = "  -3.68%  "

$ws.Range("E35").Value = "  +11.18%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "69.70"
$ws.Range("E36").Value = "  -3.04%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "632.94"
$ws.Range("E37").Value = "  -3.92%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.439"
$ws.Range("E38").Value = "  -0.44%  "

$ws.Range("E39").Value = "  +0.46%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.998"
$ws.Range("E40").Value = "  -0.10%  "

$ws.Range("E41").Value = "  -0.08%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.28"
$ws.Range("E42").Value = "  -2.53%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.87"
$ws.Range("E43").Value = "  +7.91%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.17"
$ws.Range("E44").Value = "  +19.56%  "

$ws.Range("E45").Value = "  -3.51%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "10.21"
$ws.Range("E46").Value = "  -3.98%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.94"
$ws.Range("E47").Value = "  -12.18%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.145"
$ws.Range("E48").Value = "  -3.79%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.907.77"
$ws.Range("E49").Value = "  -0.14%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "3.24"
$ws.Range("E50").Value = "  -4.44%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.000278"
$ws.Range("E51").Value = "  +3.08%  "
